# Regenerate the localization-status report for archive:
#  - flip the in-flight status text from "Ready for handoff" to "In Translation"
#    on every sheet that surfaces it (the Overview roll-up plus each locale tab)
#  - shrink the now-narrower Status/locale columns to fit the new text

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns (E2, F2) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E1:F1").ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C2) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C1").ColumnWidth = 12.5

# --- de-de sheet: Status column (C2) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C1").ColumnWidth = 12.5
